$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 10; Date = "08/01/2025"; B = 10.36499999999999;  C = 3.787;              D = 14.10400000000001; E = 0.3799999999999998;  F = 0;   G = 265.15; H = 12.206;             I = 27.40662903225803; J = 1015.845 },
    @{ Row = 11; Date = "09/01/2025"; B = 9.668999999999997;  C = 6.206;              D = 17.046;             E = 0.7440000000000001;  F = 0.8; G = 289.39; H = 15.147;             I = 27.81462903225803; J = 1018.409 },
    @{ Row = 12; Date = "10/01/2025"; B = 12.037;             C = 3.449999999999999;  D = 13.979;             E = 1.573;                F = 0;   G = 278.96; H = 16.84399999999999;  I = 27.34816129032254; J = 1016.631 }
)

# Scratch row used to stage text values (cleared at the end) so that
# assigning a dd/mm/yyyy-looking string via PasteSpecial(values-only)
# doesn't disturb the destination cell's existing format/style.
$scratchRow = 1000

foreach ($r in $rows) {
    $row = $r.Row

    # Duplicate the formatting (and values) of the last existing data row
    # (row 9) onto the new row so style/number-format stay consistent.
    $ws.Range("A9:K9").Copy($ws.Range("A" + $row + ":K" + $row))

    # Column A holds a dd/mm/yyyy-look-alike string. Stage it in a
    # text-formatted scratch cell so Excel doesn't reinterpret it as a
    # date serial, then paste only the VALUE into place so A10's style
    # (copied above from A9) is preserved untouched.
    $ws.Cells.Item($scratchRow, 1).NumberFormat = "@"
    $ws.Cells.Item($scratchRow, 1).Value = $r.Date
    $ws.Cells.Item($scratchRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = "Valencia"
}

# Drop the scratch row entirely so it doesn't affect the sheet's used
# range/dimension.
$ws.Rows.Item($scratchRow).Delete()
